# Apply cryptos list update (prices and 1h volume %) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay as text, even when they look numeric
# (e.g. "0.820", "132.80", "0.0510") so trailing zeros/format are preserved.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.168.98"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.522.52"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "536.29"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "137.91"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "2.521.08"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").Value = "5.35"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").Value = "2.956.98"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "23.03"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "59.033.96"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "2.527.90"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "11.12"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "325.36"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("D24").Value = "65.79"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").Value = "6.68"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +6.34%  "
$ws.Range("D33").Value = "163.86"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "36.61"
$ws.Range("D40").Value = "0.820"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("D42").Value = "286.45"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "132.80"
$ws.Range("E44").Value = "  +7.65%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "10.89"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -2.21%  "

# Remove the temporary text-format styling so cells keep the original (default) style
$ws.Range("D2:D51").ClearFormats()
